# add the NA's under duplicate_image_filename (column E), rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
